# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`r`n" +
           "✅ Dólar paralelo: 68`r`n`r`n" +
           "Binance`r`n" +
           "✅ 1000 Bs = 4.67 = 18443.97 pesos`r`n" +
           "✅ 18443.97 pesos = 4.65 = 964.58 Bs`r`n`r`n" +
           "Promedio competencia`r`n" +
           "✅ Tasa pesos: 20`r`n" +
           "✅ Tasa Bs: 20`r`n" +
           "✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 214
$wsTasas.Range("O10").Value = 3947.01
$wsTasas.Range("N12").Value = 3970
$wsTasas.Range("O12").Value = 207.622
